$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = "'" + $text
}

# ---- Row 15 ----
$ws.Cells.Item(15, 1).Value = 112079417
$ws.Cells.Item(15, 2).Value = 8377
$ws.Cells.Item(15, 3).Value = "Ovaliderad"
$ws.Cells.Item(15, 4).Value = "LC"
$ws.Cells.Item(15, 5).Value = 106545
$ws.Cells.Item(15, 6).Value = "Mindre märgborre"
$ws.Cells.Item(15, 7).Value = "Tomicus minor"
$ws.Cells.Item(15, 8).Value = "(Hartig, 1834)"
$ws.Cells.Item(15, 16).Value = "Köpingsåsen, Srm"
$ws.Cells.Item(15, 17).Value = 563452.2161451608
$ws.Cells.Item(15, 18).Value = 6576050.619629455
$ws.Cells.Item(15, 19).Value = 5
$ws.Cells.Item(15, 20).Value = "Södermanland"
$ws.Cells.Item(15, 21).Value = "Eskilstuna"
$ws.Cells.Item(15, 22).Value = "Södermanland"
$ws.Cells.Item(15, 23).Value = "Öja"
Set-TextCell 15 25 "2022-09-01"
$ws.Cells.Item(15, 26).Value = "00:00"
Set-TextCell 15 27 "2022-10-31"
$ws.Cells.Item(15, 28).Value = "00:00"
$ws.Cells.Item(15, 30).Value = $false
$ws.Cells.Item(15, 31).Value = $false
$ws.Cells.Item(15, 33).Value = $false
$ws.Cells.Item(15, 49).Value = "Ralf Lundmark"
$ws.Cells.Item(15, 50).Value = "Ralf Lundmark"
$ws.Cells.Item(15, 51).Value = "Lst D inventering sandbarrskogar"

# ---- Row 16 ----
$ws.Cells.Item(16, 1).Value = 112079439
$ws.Cells.Item(16, 2).Value = 90689
$ws.Cells.Item(16, 3).Value = "Ovaliderad"
$ws.Cells.Item(16, 4).Value = "NT"
$ws.Cells.Item(16, 5).Value = 5966
$ws.Cells.Item(16, 6).Value = "Motaggsvamp"
$ws.Cells.Item(16, 7).Value = "Sarcodon squamosus"
$ws.Cells.Item(16, 8).Value = "(Schaeff.) Quél."
$ws.Cells.Item(16, 16).Value = "Köpingsåsen, Srm"
$ws.Cells.Item(16, 17).Value = 563407.8844683191
$ws.Cells.Item(16, 18).Value = 6576468.999423527
$ws.Cells.Item(16, 19).Value = 5
$ws.Cells.Item(16, 20).Value = "Södermanland"
$ws.Cells.Item(16, 21).Value = "Eskilstuna"
$ws.Cells.Item(16, 22).Value = "Södermanland"
$ws.Cells.Item(16, 23).Value = "Öja"
Set-TextCell 16 25 "2022-09-01"
$ws.Cells.Item(16, 26).Value = "00:00"
Set-TextCell 16 27 "2022-10-31"
$ws.Cells.Item(16, 28).Value = "00:00"
$ws.Cells.Item(16, 30).Value = $false
$ws.Cells.Item(16, 31).Value = $false
$ws.Cells.Item(16, 33).Value = $false
$ws.Cells.Item(16, 49).Value = "Ralf Lundmark"
$ws.Cells.Item(16, 50).Value = "Ralf Lundmark"
$ws.Cells.Item(16, 51).Value = "Lst D inventering sandbarrskogar"
